## Applies "Add files via upload" edit: 8 data rows in the CVE dataset were
## replaced with different CVE records (rows 76-79 within Group 27, and rows
## 88-91 within Group 30). All other rows are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group 27 (rows 76-79) ---------------------------------------------
$ws.Range("B76").Value = 'CVE-2018-15698'
$ws.Range("C76").Value = 'ASUSTOR Data Master 3.1.5 and below allows authenticated remote non-administrative users to read any file on the file system when providing the full path to loginimage.cgi.'
$ws.Range("D76").Value = 'providing the full path to loginimage.cgi'
$ws.Range("E76").Value = 'read any file on the file system'
$ws.Range("F76").Value = 'cpe:2.3:o:asustor:data_master:*:*:*:*:*:*:*:*'
$ws.Range("G76").Value = 8

$ws.Range("B77").Value = 'CVE-2008-6278'
$ws.Range("C77").Value = 'Multiple cross-site scripting (XSS) vulnerabilities in product.php in RakhiSoftware Price Comparison Script (aka Shopping Cart) allow remote attackers to inject arbitrary web script or HTML via the (1) category_id and (2) subcategory_id parameters.'
$ws.Range("D77").Value = 'the (1) category_id and (2) subcategory_id parameters'
$ws.Range("E77").Value = 'inject arbitrary web script or HTML'
$ws.Range("F77").Value = 'cpe:2.3:a:rakhisoftware:rakhisoftware_shopping_cart:-:*:*:*:*:*:*:*'
$ws.Range("G77").Value = 8.6

$ws.Range("B78").Value = 'CVE-1999-0978'
$ws.Range("C78").Value = 'htdig allows remote attackers to execute commands via filenames with shell metacharacters.'
$ws.Range("D78").Value = 'filenames with shell metacharacters'
$ws.Range("E78").Value = 'execute commands'
$ws.Range("F78").Value = 'cpe:2.3:o:debian:debian_linux:2.1:*:*:*:*:*:*:*'
$ws.Range("G78").Value = 10

$ws.Range("B79").Value = 'CVE-2018-20823'
$ws.Range("C79").Value = 'The gyroscope on Xiaomi Mi 5s devices allows attackers to cause a denial of service (resonance and false data) via a 20.4 kHz audio signal, aka a MEMS ultrasound attack.'
$ws.Range("D79").Value = 'a 20.4 kHz audio signal'
$ws.Range("E79").Value = 'cause a denial of service (resonance and false data)'
$ws.Range("F79").Value = 'cpe:2.3:o:mi:mi_5s_firmware:-:*:*:*:*:*:*:*'
$ws.Range("G79").Value = 10

# --- Group 30 (rows 88-91) ----------------------------------------------
$ws.Range("B88").Value = 'CVE-2017-5609'
$ws.Range("C88").Value = 'SQL injection vulnerability in include/functions_entries.inc.php in Serendipity 2.0.5 allows remote authenticated users to execute arbitrary SQL commands via the cat parameter.'
$ws.Range("D88").Value = 'the cat parameter'
$ws.Range("E88").Value = 'execute arbitrary SQL commands'
$ws.Range("F88").Value = 'cpe:2.3:a:s9y:serendipity:2.0.5:*:*:*:*:*:*:*'
$ws.Range("G88").Value = 8

$ws.Range("B89").Value = 'CVE-2009-1311'
$ws.Range("C89").Value = 'Mozilla Firefox before 3.0.9 and SeaMonkey before 1.1.17 allow user-assisted remote attackers to obtain sensitive information via a web page with an embedded frame, which causes POST data from an outer page to be sent to the inner frame''s URL during a SAVEMODE_FILEONLY save of the inner frame.'
$ws.Range("D89").Value = 'a web page with an embedded frame'
$ws.Range("E89").Value = 'obtain sensitive information'
$ws.Range("F89").Value = 'cpe:2.3:a:mozilla:firefox:0.1:*:*:*:*:*:*:*'
$ws.Range("G89").Value = 8.6

$ws.Range("B90").Value = 'CVE-2015-1344'
$ws.Range("C90").Value = 'The do_write_pids function in lxcfs.c in LXCFS before 0.12 does not properly check permissions, which allows local users to gain privileges by writing a pid to the tasks file.'
$ws.Range("D90").Value = 'writing a pid to the tasks file'
$ws.Range("E90").Value = 'gain privileges'
$ws.Range("F90").Value = 'cpe:2.3:o:canonical:ubuntu_linux:15.04:*:*:*:*:*:*:*'
$ws.Range("G90").Value = 3.9

$ws.Range("B91").Value = 'CVE-2002-2273'
$ws.Range("C91").Value = 'Cross-site scripting (XSS) vulnerability in Webster HTTP Server allows remote attackers to inject arbitrary web script or HTML via the URL.'
$ws.Range("D91").Value = 'the URL'
$ws.Range("E91").Value = 'inject arbitrary web script or HTML'
$ws.Range("F91").Value = 'cpe:2.3:a:webster:webster_http_server:*:*:*:*:*:*:*:*'
$ws.Range("G91").Value = 8.6

# --- Restore the view/selection state recorded in the saved file --------
$ws.Application.ActiveWindow.ScrollRow = 78
$ws.Range("I88").Select()
